$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the header: "price_aftert_discount" -> "price_after_discount"
$ws.Range("I1").Value = "price_after_discount"

# Update the active cell/selection to I1 (matches the diff's sheetView selection)
$ws.Range("I1").Select()
